# Insert a new data row before the current row 63 (shifts rows 63..123 down
# to 64..124, same as pressing "Insert" on that row in Excel).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(63).Insert()

# Populate the newly inserted row 63 with the new market-report record.
$ws.Cells.Item(63, 1).Value = 5
$ws.Cells.Item(63, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(63, 3).Value = "Maule"
$ws.Cells.Item(63, 4).Value = 44512
$ws.Cells.Item(63, 5).Value = 7
$ws.Cells.Item(63, 6).Value = 100112024
$ws.Cells.Item(63, 7).Value = "Choclo"
$ws.Cells.Item(63, 8).Value = "Dulce o Americano"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 300
$ws.Cells.Item(63, 11).Value = 20000
$ws.Cells.Item(63, 12).Value = 20000
$ws.Cells.Item(63, 13).Value = 20000
$ws.Cells.Item(63, 14).Value = "$/malla 60 unidades"
$ws.Cells.Item(63, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(63, 16).Value = 333
$ws.Cells.Item(63, 17).Value = 60
$ws.Cells.Item(63, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(63, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
